$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Step 1: B1 currently holds the shared string "original value" and is its sole user.
# Overwriting it in-place lets that shared-string slot flip its text to "editor"
# (matching the later-reused index that the new column E header will point to).
$ws.Range("B1").Value = "editor"

# Step 2: Fill the new column E (rows 2-31) with the placeholder value "-"
$ws.Range("E2:E31").Value = "-"

# Step 3: Give the new column E a header in row 1 with the same text that B1 now has
# ("editor") - this reuses the shared string created in step 1.
$ws.Range("E1").Value = "editor"

# Copy A1's formatting (bold font + header fill) onto E1 without touching its value.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# Step 4: Now change B1 to its real header text. Since E1 also uses the "editor"
# shared string now, this creates a new shared-string entry for B1's new text.
$ws.Range("B1").Value = "Persian (TRN-1079-0833-9890)"

# Update the active selection to B2, matching the saved view state.
$ws.Range("B2").Select()
